$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# New row (2021) of statistics, appended directly below the existing data
# (row 11 holds the 2020 figures). Columns C, D, O, S, W have no reported
# value for 2021 and stay blank, just like the equivalent blank cells used
# for other years.
$rowNumber = 12

# Reuse the formatting of the previous row's cells (Copy + PasteSpecial of
# formats only) so the new row shares the exact same styles already defined
# in the workbook instead of creating new, duplicate style entries.
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A$rowNumber").PasteSpecial($xlPasteFormats) | Out-Null

$blankColumns = @("C", "D", "O", "S", "W")
foreach ($col in $blankColumns) {
    $ws.Range("$col" + "11").Copy() | Out-Null
    $ws.Range("$col$rowNumber").PasteSpecial($xlPasteFormats) | Out-Null
}

$excel.CutCopyMode = 0

$ws.Range("A$rowNumber").Value = "2021年"

$numericValues = @{
    "B"  = 207
    "E"  = 7244
    "F"  = 69767
    "G"  = 23
    "H"  = 480020
    "I"  = 63
    "J"  = 9187
    "K"  = 127803
    "L"  = 392
    "M"  = 157
    "N"  = 29503
    "P"  = 783
    "Q"  = 1226165
    "R"  = 41
    "T"  = 9611
    "U"  = 50403
    "V"  = 639548
    "X"  = 60154
    "Y"  = 13
    "Z"  = 403
    "AA" = 5
    "AB" = 552
    "AC" = 1
}

foreach ($col in $numericValues.Keys) {
    $ws.Range("$col$rowNumber").Value = $numericValues[$col]
}
